$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "MaxAnalysisLevel"
$ws.Range("B17").Value = "anlMax"
$ws.Range("D17").Value = 110
